# Daily attendance processing - 2025-12-05 13:40:28
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# stores a comma-separated list of the users who recorded/edited a
# session's attendance. This pass normalizes the ordering of that list
# for the specific pairings that were recorded in the wrong order:
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "System, admin@admin.com"             -> "admin@admin.com, System"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
# Any other "Recorded By" values (single names, three-way combinations,
# or pairs not listed above) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$swapMap = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # Column G = "Recorded By"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Value2

    if ($null -ne $current -and $swapMap.ContainsKey($current)) {
        $cell.Value = $swapMap[$current]
    }
}
